# Recommandations sheet: refresh BRVM sector-index rows + top-mover stock rows
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$data1 = @(
    @(2, "BRVM - CONSOMMATION DE BASE     (**)", 0, 4, 903.43, 231.36, "🟡 Observer", "➖ Neutre"),
    @(3, "BRVM-PRINCIPAL     (**)", 0, 4, 902.91, 228.7, "🟡 Observer", "➖ Neutre"),
    @(4, "BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 681.65, 171.76, "🟡 Observer", "➖ Neutre"),
    @(5, "BRVM - SERVICES FINANCIERS", 0, 4, 592.28, 149.63, "🟡 Observer", "➖ Neutre"),
    @(6, "BRVM - INDUSTRIELS", 0, 4, 578.55, 146.16, "🟡 Observer", "➖ Neutre"),
    @(7, "BRVM-PRESTIGE", 0, 4, 575.83, 145.37, "🟡 Observer", "➖ Neutre"),
    @(8, "BRVM – COMPOSITE TOTAL RETURN     (**)", 0, 4, 538.69, 136, "🟡 Observer", "➖ Neutre"),
    @(9, "BRVM - ENERGIE", 0, 4, 455.75, 114.43, "🟡 Observer", "➖ Neutre"),
    @(10, "BRVM - SERVICES PUBLICS", 0, 4, 453.06, 114.45, "🟡 Observer", "➖ Neutre"),
    @(11, "BRVM - TELECOMMUNICATIONS", 0, 4, 376.07, 94.38, "🟡 Observer", "➖ Neutre"),
    @(12, "EVIOSYS PACKAGING SIEM CI (SEMC)", 4, 0, 24.54, 6.49, "🟢 Achat", "✅ Renforcer"),
    @(13, "UNILEVER CI (UNLC)", 2, 0, 14.32, 6.83, "🟡 Observer", "➖ Neutre"),
    @(14, "CFAO MOTORS CI (CFAC)", 2, 1, 7.72, -0.69, "🟡 Observer", "👀 À surveiller"),
    @(15, "TRACTAFRIC MOTORS CI (PRSC)", 1, 0, 5.25, 5.25, "🟡 Observer", "➖ Neutre"),
    @(16, "NSIA BANQUE COTE D'IVOIRE (NSBC)", 1, 0, 5.04, 5.04, "🟡 Observer", "➖ Neutre"),
    @(17, "BERNABE CI (BNBC)", 1, 1, 4.88, 7.41, "🟡 Observer", "👀 À surveiller"),
    @(18, "SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 1, 0, 3.35, 3.35, "🟡 Observer", "➖ Neutre"),
    @(19, "CORIS BANK INTERNATIONAL (CBIBF)", 1, 1, 3.22, 7.11, "🟡 Observer", "👀 À surveiller"),
    @(20, "SICOR CI (SICC)", 1, 0, 3.11, 3.11, "🟡 Observer", "➖ Neutre"),
    @(21, "NESTLE CI (NTLC)", 2, 1, 2.27, 3.64, "🟡 Observer", "👀 À surveiller"),
    @(22, "SOLIBRA CI (SLBC)", 1, 1, 0.36, 6.08, "🟡 Observer", "👀 À surveiller"),
    @(23, "AFRICA GLOBAL LOGISTICS CI (SDSC)", 1, 2, 0.26, -2.89, "🟡 Observer", "👀 À surveiller"),
    @(24, "BICI CI (BICC)", 0, 1, -0.49, -0.49, "🟡 Observer", "➖ Neutre"),
    @(25, "SICABLE CI (CABC)", 1, 1, -0.68, -3.35, "🟡 Observer", "👀 À surveiller"),
    @(26, "FILTISAC CI (FTSC)", 0, 1, -0.7, -0.7, "🟡 Observer", "➖ Neutre"),
    @(27, "TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -1.88, -1.88, "🟡 Observer", "➖ Neutre"),
    @(28, "TOTALENERGIES MARKETING SN (TTLS)", 0, 1, -1.92, -1.92, "🟡 Observer", "➖ Neutre"),
    @(29, "ORAGROUP TOGO (ORGT)", 0, 1, -2.08, -2.08, "🟡 Observer", "➖ Neutre"),
    @(30, "LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -2.44, -2.44, "🟡 Observer", "➖ Neutre"),
    @(31, "SMB CI (SMBC)", 0, 1, -2.57, -2.57, "🟡 Observer", "➖ Neutre"),
    @(32, "BANK OF AFRICA BF (BOABF)", 0, 1, -3.49, -3.49, "🟡 Observer", "➖ Neutre"),
    @(33, "ECOBANK TRANS. INCORP. TG (ETIT)", 1, 2, -4.15, -4.35, "🟡 Observer", "👀 À surveiller"),
    @(34, "NEI-CEDA CI (NEIC)", 0, 1, -4.26, -4.26, "🟡 Observer", "➖ Neutre"),
    @(35, "SETAO CI (STAC)", 0, 1, -4.35, -4.35, "🟡 Observer", "➖ Neutre")
)

foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
    $ws1.Cells.Item($r, 7).Value = $row[7]
}

$ws1.Range("A1:G35").EntireColumn.AutoFit() | Out-Null

# Top_YTD sheet: refresh BRVM sector-index YTD progression table
$ws2 = $wb.Worksheets.Item("Top_YTD")

$data2 = @(
    @(2, "BRVM - CONSOMMATION DE BASE     (**)", 11170.88),
    @(3, "BRVM-PRINCIPAL     (**)", 11155.55),
    @(4, "BRVM - CONSOMMATION DISCRETIONNAIRE", 5246.3),
    @(5, "BRVM - SERVICES FINANCIERS", 3686.88),
    @(6, "BRVM - INDUSTRIELS", 3481.06),
    @(7, "BRVM-PRESTIGE", 3441.97),
    @(8, "BRVM – COMPOSITE TOTAL RETURN     (**)", 2932.72),
    @(9, "BRVM - ENERGIE", 1994.8),
    @(10, "BRVM - SERVICES PUBLICS", 1968.54),
    @(11, "BRVM - TELECOMMUNICATIONS", 1316.96)
)

foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
}
